$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - same formatting as the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cells H2/H3 with value 0 (plain, unstyled, like the rest of the data column)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
